$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 392.66666
$ws.Range("I4").Value = 140.25
$ws.Range("K4").Value = 140.25
$ws.Range("M4").Value = -26.25
$ws.Range("H33").Value = 87.666664
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H58").Value = 1287.7
$ws.Range("J58").Value = 1792.8572
$ws.Range("L58").Value = 5378.571599999999
$ws.Range("N58").Value = -5678.571599999999
$ws.Range("H69").Value = 6882.5947
$ws.Range("I69").Value = 9995
$ws.Range("J69").Value = 6704.7427
$ws.Range("K69").Value = 29985
$ws.Range("L69").Value = 20114.2281
$ws.Range("M69").Value = -29111
$ws.Range("N69").Value = -21862.2281
$ws.Range("H72").Value = 6882.5947
$ws.Range("I72").Value = 9995
$ws.Range("J72").Value = 6704.7427
$ws.Range("K72").Value = 89955
$ws.Range("L72").Value = 60342.6843
$ws.Range("M72").Value = -85587
$ws.Range("N72").Value = -69078.68429999999
$ws.Range("H100").Value = 3022.5557
$ws.Range("I100").Value = 533.8333
$ws.Range("J100").Value = 8000
$ws.Range("K100").Value = 533.8333
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = 7.166699999999992
$ws.Range("N100").Value = -9082
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6968.92
$ws.Range("I32").Value = 4601.091
$ws.Range("J32").Value = 24333
$ws.Range("K32").Value = 4601.091
$ws.Range("L32").Value = 24333
$ws.Range("M32").Value = -4314.091
$ws.Range("N32").Value = -24907
$ws.Range("H45").Value = 3552.6667
$ws.Range("I45").Value = 2495.6667
$ws.Range("K45").Value = 2495.6667
$ws.Range("M45").Value = -2118.6667
$ws.Range("H74").Value = 8999.6
$ws.Range("J74").Value = 9110.666999999999
$ws.Range("L74").Value = 9110.666999999999
$ws.Range("N74").Value = -10858.667
$ws.Range("H77").Value = 8999.6
$ws.Range("J77").Value = 9110.666999999999
$ws.Range("L77").Value = 45553.335
$ws.Range("N77").Value = -54289.335
$ws.Range("H97").Value = 1121.5714
$ws.Range("I97").Value = 598.36365
$ws.Range("K97").Value = 598.36365
$ws.Range("M97").Value = -102.36365
$ws.Range("H122").Value = 1652
$ws.Range("I122").Value = 1536.5
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 4609.5
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -2159.5
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 3282.3
$ws.Range("I132").Value = 791
$ws.Range("J132").Value = 5773.6
$ws.Range("K132").Value = 2373
$ws.Range("L132").Value = 17320.8
$ws.Range("M132").Value = 157
$ws.Range("N132").Value = -22380.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 1492
$ws.Range("J15").Value = 1500
$ws.Range("L15").Value = 1500
$ws.Range("N15").Value = -1954
$ws.Range("H107").Value = 5249.095
$ws.Range("I107").Value = 1733.3
$ws.Range("K107").Value = 1733.3
$ws.Range("M107").Value = 186.7
$ws.Range("H134").Value = 3263
$ws.Range("I134").Value = 2758.1333
$ws.Range("J134").Value = 7049.5
$ws.Range("K134").Value = 8274.3999
$ws.Range("L134").Value = 21148.5
$ws.Range("M134").Value = -5739.3999
$ws.Range("N134").Value = -26218.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2673.75
$ws.Range("I86").Value = 2481.6667
$ws.Range("K86").Value = 2481.6667
$ws.Range("M86").Value = -1358.6667
$ws.Range("H89").Value = 2673.75
$ws.Range("I89").Value = 2481.6667
$ws.Range("K89").Value = 12408.3335
$ws.Range("M89").Value = -6792.333500000001
$ws.Range("H93").Value = 1500
$ws.Range("I93").Value = 1500
$ws.Range("K93").Value = 1500
$ws.Range("M93").Value = 372
$ws.Range("H99").Value = 4386.375
$ws.Range("I99").Value = 3791.3076
$ws.Range("K99").Value = 3791.3076
$ws.Range("M99").Value = -2293.3076
$ws.Range("H126").Value = 4386.375
$ws.Range("I126").Value = 3791.3076
$ws.Range("K126").Value = 11373.9228
$ws.Range("M126").Value = -8903.9228
$ws.Range("H132").Value = 5063.16
$ws.Range("I132").Value = 4318.3125
$ws.Range("J132").Value = 6387.3335
$ws.Range("K132").Value = 12954.9375
$ws.Range("L132").Value = 19162.0005
$ws.Range("M132").Value = -10424.9375
$ws.Range("N132").Value = -24222.0005
$ws.Range("H134").Value = 3688.5
$ws.Range("I134").Value = 3235.75
$ws.Range("J134").Value = 5499.5
$ws.Range("K134").Value = 9707.25
$ws.Range("L134").Value = 16498.5
$ws.Range("M134").Value = -7172.25
$ws.Range("N134").Value = -21568.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 157.61539
$ws.Range("J50").Value = 166.58333
$ws.Range("L50").Value = 499.74999
$ws.Range("N50").Value = -1461.74999
$ws.Range("H53").Value = 157.61539
$ws.Range("J53").Value = 166.58333
$ws.Range("L53").Value = 499.74999
$ws.Range("N53").Value = -1461.74999
$ws.Range("H116").Value = 1615.1111
$ws.Range("I116").Value = 1674.6666
$ws.Range("K116").Value = 5023.9998
$ws.Range("M116").Value = -1581.9998
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 23333
$ws.Range("I63").Value = 23333
$ws.Range("K63").Value = 23333
$ws.Range("M63").Value = -22647
$ws.Range("H66").Value = 23333
$ws.Range("I66").Value = 23333
$ws.Range("K66").Value = 69999
$ws.Range("M66").Value = -66567
$ws.Range("H80").Value = 3116.25
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3155
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3155
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5151
$ws.Range("H83").Value = 3116.25
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3155
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15775
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -25759
$ws.Range("H97").Value = 620.4286
$ws.Range("I97").Value = 620.4286
$ws.Range("K97").Value = 620.4286
$ws.Range("M97").Value = -124.4286
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 500
$ws.Range("N107").Value = -4340
$ws.Range("H132").Value = 5460.375
$ws.Range("J132").Value = 8993
$ws.Range("L132").Value = 26979
$ws.Range("N132").Value = -32039
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7999.5
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10271
$ws.Range("H46").Value = 6182.2856
$ws.Range("I46").Value = 5704.6665
$ws.Range("J46").Value = 6540.5
$ws.Range("K46").Value = 5704.6665
$ws.Range("L46").Value = 6540.5
$ws.Range("M46").Value = -5516.6665
$ws.Range("N46").Value = -6916.5
$ws.Range("H82").Value = 4393.533
$ws.Range("J82").Value = 6575.375
$ws.Range("L82").Value = 6575.375
$ws.Range("N82").Value = -7297.375
$ws.Range("H85").Value = 4393.533
$ws.Range("J85").Value = 6575.375
$ws.Range("L85").Value = 6575.375
$ws.Range("N85").Value = -9071.375
$ws.Range("H93").Value = 1188.0667
$ws.Range("I93").Value = 1220.9166
$ws.Range("J93").Value = 1056.6666
$ws.Range("K93").Value = 1220.9166
$ws.Range("L93").Value = 1056.6666
$ws.Range("M93").Value = 27.08339999999998
$ws.Range("N93").Value = -3552.6666
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("N129").Value = 0
$ws.Range("H141").Value = 206125
$ws.Range("J141").Value = 206125
$ws.Range("L141").Value = 206125
$ws.Range("N141").Value = -216485
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 203406
$ws.Range("I2").Value = 203406
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 203406
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -203294
$ws.Range("H64").Value = 49989
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 49989
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H98").Value = 38300
$ws.Range("J98").Value = 38300
$ws.Range("L98").Value = 38300
$ws.Range("N98").Value = -44290
$ws.Range("H107").Value = 473.42105
$ws.Range("I107").Value = 406.7857
$ws.Range("K107").Value = 1220.3571
$ws.Range("M107").Value = 699.6428999999998
$ws.Range("H132").Value = 2559.2144
$ws.Range("I132").Value = 1923.25
$ws.Range("K132").Value = 5769.75
$ws.Range("M132").Value = -3239.75
$ws.Range("H136").Value = 3377.45
$ws.Range("I136").Value = 2141.3333
$ws.Range("J136").Value = 4388.8184
$ws.Range("K136").Value = 6423.999899999999
$ws.Range("L136").Value = 13166.4552
$ws.Range("M136").Value = -3873.999899999999
$ws.Range("N136").Value = -18266.4552
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280